$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows at 222-223 (existing rows 222-250 shift down to 224-252)
$ws.Rows("222:223").Insert()

# New row 222: Terminal Hortofrutícola Agro Chillán - Repollo, Crespo record, Primera
$ws.Cells.Item(222,1).Value = 7
$ws.Cells.Item(222,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(222,3).Value = "Ñuble"
$ws.Cells.Item(222,4).Value = 44776
$ws.Cells.Item(222,5).Value = 16
$ws.Cells.Item(222,6).Value = 100112006
$ws.Cells.Item(222,7).Value = "Repollo"
$ws.Cells.Item(222,8).Value = "Crespo record"
$ws.Cells.Item(222,9).Value = "Primera"
$ws.Cells.Item(222,10).Value = 200
$ws.Cells.Item(222,11).Value = 1000
$ws.Cells.Item(222,12).Value = 1200
$ws.Cells.Item(222,13).Value = 1100
$ws.Cells.Item(222,14).Value = "$/unidad"
$ws.Cells.Item(222,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(222,16).Value = 1100
$ws.Cells.Item(222,17).Value = 1
$ws.Cells.Item(222,18).Value = "Hortaliza"

# New row 223: Terminal Hortofrutícola Agro Chillán - Repollo, Crespo record, Segunda
$ws.Cells.Item(223,1).Value = 7
$ws.Cells.Item(223,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(223,3).Value = "Ñuble"
$ws.Cells.Item(223,4).Value = 44776
$ws.Cells.Item(223,5).Value = 16
$ws.Cells.Item(223,6).Value = 100112006
$ws.Cells.Item(223,7).Value = "Repollo"
$ws.Cells.Item(223,8).Value = "Crespo record"
$ws.Cells.Item(223,9).Value = "Segunda"
$ws.Cells.Item(223,10).Value = 150
$ws.Cells.Item(223,11).Value = 900
$ws.Cells.Item(223,12).Value = 900
$ws.Cells.Item(223,13).Value = 900
$ws.Cells.Item(223,14).Value = "$/unidad"
$ws.Cells.Item(223,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(223,16).Value = 900
$ws.Cells.Item(223,17).Value = 1
$ws.Cells.Item(223,18).Value = "Hortaliza"
